$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume columns) with the latest scrape.
# Column D holds price strings that must stay TEXT (some use '.' as a
# thousands separator, e.g. "29.039.06", others are small decimals like
# "0.000009339"); a bare numeric literal would get auto-coerced by Excel
# into a real number (dropping trailing zeros / flipping to exponential
# notation), so those assignments use a leading apostrophe quote-prefix,
# same as a user typing '0.000009339 into the cell, to force text entry.
$ws.Range("D2").Value = '''29.039.06'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '''1.830.10'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''241.15'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").Value = '''0.6263'
$ws.Range("E6").Value = '  -5.18%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '''45.22'
$ws.Range("E8").Value = '  +8.12%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '''0.07613'
$ws.Range("E9").Value = '  +2.75%  '
$ws.Range("D10").Value = '''0.2913'
$ws.Range("E10").Value = '  -0.78%  '
$ws.Range("D11").Value = '''22.75'
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").Value = '''0.07645'
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("D13").Value = '''1.829.51'
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").Value = '''4.959'
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").Value = '''82.34'
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("D17").Value = '''0.000009339'
$ws.Range("E17").Value = '  +9.74%  '
$ws.Range("D18").Value = '''5.991'
$ws.Range("D19").Value = '''28.863.36'
$ws.Range("D20").Value = '''224.79'
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("D21").Value = '''12.32'
$ws.Range("E21").Value = '  -1.05%  '
$ws.Range("D22").Value = '''0.9998'
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").Value = '''7.207'
$ws.Range("E23").Value = '  +1.79%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").Value = '''159.52'
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = '''8.419'
$ws.Range("E26").Value = '  -2.15%  '
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("D28").Value = '''17.81'
$ws.Range("D29").Value = '''1.493'
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("E30").Value = '  -1.51%  '
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("E32").Value = '  +0.83%  '
$ws.Range("D33").Value = '''0.05198'
$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("E34").Value = '  -0.84%  '
$ws.Range("D35").Value = '''1.154'
$ws.Range("D36").Value = '''0.7308'
$ws.Range("E36").Value = '  -0.62%  '
$ws.Range("D37").Value = '''2.612'
$ws.Range("E37").Value = '  -1.81%  '
$ws.Range("D38").Value = '''1.275.96'
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").Value = '''0.01787'
$ws.Range("E40").Value = '  -0.34%  '
$ws.Range("D41").Value = '''6.529'
$ws.Range("E41").Value = '  +7.75%  '
$ws.Range("D42").Value = '''0.8902'
$ws.Range("E42").Value = '  -3.32%  '
$ws.Range("D43").Value = '''0.9997'
$ws.Range("D44").Value = '''101.46'
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("D45").Value = '''1.975.84'
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("D46").Value = '''0.5106'
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("D47").Value = '''63.67'
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").Value = '''0.3977'
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").Value = '''0.07306'
$ws.Range("E50").Value = '  -13.40%  '
$ws.Range("D51").Value = '''8.861'
$ws.Range("E51").Value = '  +1.44%  '
